# Sheet1 gains a new "备注" (remark) column G: whether each API could be
# implemented ("实现" / "不能实现"), plus F18's fail-state text is rewritten.
# The view also scrolled down a bit and the selection moved to G14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- new column G: width + per-row "implemented?" note ------------------
# xlCenter = -4108 ; matches the existing vertical-center style (s="1")
# already used throughout column A, so re-use it rather than minting a
# new cell style.
$ws.Columns.Item(7).ColumnWidth = 15.78

$gValues = @{
    2  = "不能实现"
    3  = "实现"
    4  = "实现"
    5  = "实现"
    6  = "不能实现"
    7  = "实现"
    8  = "实现"
    9  = "实现"
    10 = "实现"
    11 = "实现"
    12 = "实现"
    13 = "实现"
    14 = "不能实现"
    15 = "实现"
    16 = "实现"
    17 = "实现"
    18 = "实现"
}

foreach ($row in $gValues.Keys) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.VerticalAlignment = -4108
    $cell.Value = $gValues[$row]
}

# --- F18 fail-state text rewritten --------------------------------------
$ws.Range("F18").Value = "String state:状态i=2[fail]"

# --- view: scrolled to row 10, selection now on G14 ---------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select()
